# Splits the two run-long "Programa" paragraphs (PT + EN/italic) into
# sentence-per-line blocks by inserting manual line breaks (^l -> <w:br/>)
# at each sentence boundary, matching the target OOXML diff.
$d = $word.ActiveDocument
$failures = New-Object System.Collections.ArrayList

# Portuguese "Programa" paragraph: insert line breaks between sentences
$ok = $d.Content.Find.Execute("E Arduino.Conceitos básic", $false, $false, $false, $false, $false, $true, 1, $false, "E Arduino.^lConceitos básic", 2)
if (-not $ok) { [void]$failures.Add("E Arduino.Conceitos básic") }
$ok = $d.Content.Find.Execute("etrônicos.Introdução à li", $false, $false, $false, $false, $false, $true, 1, $false, "etrônicos.^lIntrodução à li", 2)
if (-not $ok) { [void]$failures.Add("etrônicos.Introdução à li") }
$ok = $d.Content.Find.Execute("ibliotecasEntradas e saíd", $false, $false, $false, $false, $false, $true, 1, $false, "ibliotecas^lEntradas e saíd", 2)
if (-not $ok) { [void]$failures.Add("ibliotecasEntradas e saíd") }
$ok = $d.Content.Find.Execute(" digitais.Controle de dis", $false, $false, $false, $false, $false, $true, 1, $false, " digitais.^lControle de dis", 2)
if (-not $ok) { [void]$failures.Add(" digitais.Controle de dis") }
$ok = $d.Content.Find.Execute("zando PWM.Eletrônica anal", $false, $false, $false, $false, $false, $true, 1, $false, "zando PWM.^lEletrônica anal", 2)
if (-not $ok) { [void]$failures.Add("zando PWM.Eletrônica anal") }
$ok = $d.Content.Find.Execute(" Arduino. Oficina: leitur", $false, $false, $false, $false, $false, $true, 1, $false, " Arduino. ^lOficina: leitur", 2)
if (-not $ok) { [void]$failures.Add(" Arduino. Oficina: leitur") }
$ok = $d.Content.Find.Execute("al da IDE.Controle de mot", $false, $false, $false, $false, $false, $true, 1, $false, "al da IDE.^lControle de mot", 2)
if (-not $ok) { [void]$failures.Add("al da IDE.Controle de mot") }
$ok = $d.Content.Find.Execute("elé e SSR.Tópicos avançad", $false, $false, $false, $false, $false, $true, 1, $false, "elé e SSR.^lTópicos avançad", 2)
if (-not $ok) { [void]$failures.Add("elé e SSR.Tópicos avançad") }
$ok = $d.Content.Find.Execute("Bluetooth.Armazenamento d", $false, $false, $false, $false, $false, $true, 1, $false, "Bluetooth.^lArmazenamento d", 2)
if (-not $ok) { [void]$failures.Add("Bluetooth.Armazenamento d") }
$ok = $d.Content.Find.Execute("emória SD.Desenvolvimento", $false, $false, $false, $false, $false, $true, 1, $false, "emória SD.^lDesenvolvimento", 2)
if (-not $ok) { [void]$failures.Add("emória SD.Desenvolvimento") }
$ok = $d.Content.Find.Execute("qualidade.Desenvolvimento", $false, $false, $false, $false, $false, $true, 1, $false, "qualidade.^lDesenvolvimento", 2)
if (-not $ok) { [void]$failures.Add("qualidade.Desenvolvimento") }

# English "Programa" paragraph (italic): insert line breaks between sentences
$ok = $d.Content.Find.Execute("duino IDE.Basic concepts ", $false, $false, $false, $false, $false, $true, 1, $false, "duino IDE.^lBasic concepts ", 2)
if (-not $ok) { [void]$failures.Add("duino IDE.Basic concepts ") }
$ok = $d.Content.Find.Execute(" circuits.Introduction to", $false, $false, $false, $false, $false, $true, 1, $false, " circuits.^lIntroduction to", 2)
if (-not $ok) { [void]$failures.Add(" circuits.Introduction to") }
$ok = $d.Content.Find.Execute(" librariesArduino Inputs ", $false, $false, $false, $false, $false, $true, 1, $false, " libraries^lArduino Inputs ", 2)
if (-not $ok) { [void]$failures.Add(" librariesArduino Inputs ") }
$ok = $d.Content.Find.Execute("l signals.Device control ", $false, $false, $false, $false, $false, $true, 1, $false, "l signals.^lDevice control ", 2)
if (-not $ok) { [void]$failures.Add("l signals.Device control ") }
$ok = $d.Content.Find.Execute("using PWM.Analog electron", $false, $false, $false, $false, $false, $true, 1, $false, "using PWM.^lAnalog electron", 2)
if (-not $ok) { [void]$failures.Add("using PWM.Analog electron") }
$ok = $d.Content.Find.Execute("onverters.Workshop: Readi", $false, $false, $false, $false, $false, $true, 1, $false, "onverters.^lWorkshop: Readi", 2)
if (-not $ok) { [void]$failures.Add("onverters.Workshop: Readi") }
$ok = $d.Content.Find.Execute("l Monitor.DC motor and se", $false, $false, $false, $false, $false, $true, 1, $false, "l Monitor.^lDC motor and se", 2)
if (-not $ok) { [void]$failures.Add("l Monitor.DC motor and se") }
$ok = $d.Content.Find.Execute("y and SSR.Advanced topics", $false, $false, $false, $false, $false, $true, 1, $false, "y and SSR.^lAdvanced topics", 2)
if (-not $ok) { [void]$failures.Add("y and SSR.Advanced topics") }
$ok = $d.Content.Find.Execute("Bluetooth.Data storage us", $false, $false, $false, $false, $false, $true, 1, $false, "Bluetooth.^lData storage us", 2)
if (-not $ok) { [void]$failures.Add("Bluetooth.Data storage us") }
$ok = $d.Content.Find.Execute("mory card.Quality softwar", $false, $false, $false, $false, $false, $true, 1, $false, "mory card.^lQuality softwar", 2)
if (-not $ok) { [void]$failures.Add("mory card.Quality softwar") }
$ok = $d.Content.Find.Execute("velopment.Development of ", $false, $false, $false, $false, $false, $true, 1, $false, "velopment.^lDevelopment of ", 2)
if (-not $ok) { [void]$failures.Add("velopment.Development of ") }

if ($failures.Count -gt 0) {
    Write-Output ("FAILED anchors: " + ($failures -join " | "))
} else {
    Write-Output "All sentence breaks inserted successfully."
}
